$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.276772
$ws.Range("H2").Value = 12.830316
$ws.Range("I2").Value = 0.06135676581847978
$ws.Range("J2").Value = 0.06135676581847978
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.07074633333333
$ws.Range("N2").Value = 102.212239
$ws.Range("O2").Value = 0.5171464495142372
$ws.Range("P2").Value = 0.5171464495142373
$ws.Range("Q2").Value = 145.7128139375026
$ws.Range("R2").Value = 1311.415325437524
$ws.Range("S2").Value = 0.03173043359670333
$ws.Range("T2").Value = 0.03173043359670333

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.276772
$ws.Range("H3").Value = 12.830316
$ws.Range("I3").Value = 0.06135676581847978
$ws.Range("J3").Value = 0.06135676581847978
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.685497
$ws.Range("N3").Value = 83.056491
$ws.Range("O3").Value = 0.420227262899125
$ws.Range("P3").Value = 0.4202272628991251
$ws.Range("Q3").Value = 118.404558375684
$ws.Range("R3").Value = 1065.641025381156
$ws.Range("S3").Value = 0.02578378576024235
$ws.Range("T3").Value = 0.02578378576024236

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.276772
$ws.Range("H4").Value = 12.830316
$ws.Range("I4").Value = 0.06135676581847978
$ws.Range("J4").Value = 0.06135676581847978
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.125957666666666
$ws.Range("N4").Value = 12.377873
$ws.Range("O4").Value = 0.06262628758663766
$ws.Range("P4").Value = 0.06262628758663766
$ws.Range("Q4").Value = 17.64578022198533
$ws.Range("R4").Value = 158.812021997868
$ws.Range("S4").Value = 0.003842546461534094
$ws.Range("T4").Value = 0.003842546461534094

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 36.44531133333333
$ws.Range("H5").Value = 109.335934
$ws.Range("I5").Value = 0.5228631389891535
$ws.Range("J5").Value = 0.5228631389891535
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.07074633333333
$ws.Range("N5").Value = 102.212239
$ws.Range("O5").Value = 0.5171464495142372
$ws.Range("P5").Value = 0.5171464495142373
$ws.Range("Q5").Value = 1241.718957477358
$ws.Range("R5").Value = 11175.47061729622
$ws.Range("S5").Value = 0.2703968159101098
$ws.Range("T5").Value = 0.2703968159101099

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 36.44531133333333
$ws.Range("H6").Value = 109.335934
$ws.Range("I6").Value = 0.5228631389891535
$ws.Range("J6").Value = 0.5228631389891535
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.685497
$ws.Range("N6").Value = 83.056491
$ws.Range("O6").Value = 0.420227262899125
$ws.Range("P6").Value = 0.4202272628991251
$ws.Range("Q6").Value = 1009.006557583066
$ws.Range("R6").Value = 9081.059018247592
$ws.Range("S6").Value = 0.2197213457682567
$ws.Range("T6").Value = 0.2197213457682568

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 36.44531133333333
$ws.Range("H7").Value = 109.335934
$ws.Range("I7").Value = 0.5228631389891535
$ws.Range("J7").Value = 0.5228631389891535
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.125957666666666
$ws.Range("N7").Value = 12.377873
$ws.Range("O7").Value = 0.06262628758663766
$ws.Range("P7").Value = 0.06262628758663766
$ws.Range("Q7").Value = 150.3718117098202
$ws.Range("R7").Value = 1353.346305388382
$ws.Range("S7").Value = 0.03274497731078682
$ws.Range("T7").Value = 0.03274497731078682

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.98126466666666
$ws.Range("H8").Value = 86.943794
$ws.Range("I8").Value = 0.4157800951923667
$ws.Range("J8").Value = 0.4157800951923668
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.07074633333333
$ws.Range("N8").Value = 102.212239
$ws.Range("O8").Value = 0.5171464495142372
$ws.Range("P8").Value = 0.5171464495142373
$ws.Range("Q8").Value = 987.4133168771959
$ws.Range("R8").Value = 8886.719851894764
$ws.Range("S8").Value = 0.215019200007424
$ws.Range("T8").Value = 0.2150192000074241

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.98126466666666
$ws.Range("H9").Value = 86.943794
$ws.Range("I9").Value = 0.4157800951923667
$ws.Range("J9").Value = 0.4157800951923668
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.685497
$ws.Range("N9").Value = 83.056491
$ws.Range("O9").Value = 0.420227262899125
$ws.Range("P9").Value = 0.4202272628991251
$ws.Range("Q9").Value = 802.3607159852058
$ws.Range("R9").Value = 7221.246443866853
$ws.Range("S9").Value = 0.1747221313706259
$ws.Range("T9").Value = 0.174722131370626

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.98126466666666
$ws.Range("H10").Value = 86.943794
$ws.Range("I10").Value = 0.4157800951923667
$ws.Range("J10").Value = 0.4157800951923668
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.125957666666666
$ws.Range("N10").Value = 12.377873
$ws.Range("O10").Value = 0.06262628758663766
$ws.Range("P10").Value = 0.06262628758663766
$ws.Range("Q10").Value = 119.5754711411291
$ws.Range("R10").Value = 1076.179240270162
$ws.Range("S10").Value = 0.02603876381431674
$ws.Range("T10").Value = 0.02603876381431675
